$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 16666950
$ws.Range("I33").Value = 20000280
$ws.Range("K33").Value = 20000280
$ws.Range("M33").Value = -20000051
$ws.Range("H40").Value = 1445.1111
$ws.Range("J40").Value = 1528
$ws.Range("L40").Value = 1528
$ws.Range("N40").Value = -1878
$ws.Range("H42").Value = 150.8
$ws.Range("I42").Value = 20.75
$ws.Range("K42").Value = 62.25
$ws.Range("M42").Value = 167.75
$ws.Range("H49").Value = 214.57143
$ws.Range("I49").Value = 217
$ws.Range("J49").Value = 200
$ws.Range("K49").Value = 651
$ws.Range("L49").Value = 600
$ws.Range("M49").Value = -515
$ws.Range("N49").Value = -872
$ws.Range("H59").Value = 1200
$ws.Range("J59").Value = 1200
$ws.Range("L59").Value = 3600
$ws.Range("N59").Value = -4714
$ws.Range("H96").Value = 634.3889
$ws.Range("I96").Value = 677.8461
$ws.Range("J96").Value = 521.4
$ws.Range("K96").Value = 2033.5383
$ws.Range("L96").Value = 1564.2
$ws.Range("M96").Value = -660.5382999999999
$ws.Range("N96").Value = -4310.2
$ws.Range("H112").Value = 60540.47
$ws.Range("J112").Value = 68542.92999999999
$ws.Range("L112").Value = 205628.79
$ws.Range("N112").Value = -207844.79
$ws.Range("H132").Value = 18575504
$ws.Range("I132").Value = 20834762
$ws.Range("K132").Value = 62504286
$ws.Range("M132").Value = -62501756
$ws.Range("H135").Value = 911.0454999999999
$ws.Range("I135").Value = 419.22223
$ws.Range("K135").Value = 3773.00007
$ws.Range("M135").Value = -1238.00007
$ws.Range("H137").Value = 1461.75
$ws.Range("I137").Value = 1365.4517
$ws.Range("K137").Value = 4096.355100000001
$ws.Range("M137").Value = -1546.355100000001
$ws.Range("H138").Value = 3205.8394
$ws.Range("I138").Value = 1348.3846
$ws.Range("J138").Value = 4815.6333
$ws.Range("K138").Value = 4045.1538
$ws.Range("L138").Value = 14446.8999
$ws.Range("M138").Value = 1094.8462
$ws.Range("N138").Value = -24726.8999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6649.7393
$ws.Range("I45").Value = 7462.1113
$ws.Range("K45").Value = 7462.1113
$ws.Range("M45").Value = -7085.1113
$ws.Range("H61").Value = 3912.375
$ws.Range("I61").Value = 3559.5264
$ws.Range("J61").Value = 5253.2
$ws.Range("K61").Value = 3559.5264
$ws.Range("L61").Value = 5253.2
$ws.Range("M61").Value = -3347.5264
$ws.Range("N61").Value = -5677.2
$ws.Range("H97").Value = 2392.55
$ws.Range("I97").Value = 2638.0588
$ws.Range("J97").Value = 1001.3333
$ws.Range("K97").Value = 2638.0588
$ws.Range("L97").Value = 1001.3333
$ws.Range("M97").Value = -2142.0588
$ws.Range("N97").Value = -1993.3333
$ws.Range("H122").Value = 1324
$ws.Range("I122").Value = 1221.5238
$ws.Range("K122").Value = 3664.5714
$ws.Range("M122").Value = -1214.5714
$ws.Range("H132").Value = 4226.913
$ws.Range("I132").Value = 1891.7333
$ws.Range("K132").Value = 5675.199900000001
$ws.Range("M132").Value = -3145.199900000001
$ws.Range("H136").Value = 3912.375
$ws.Range("I136").Value = 3559.5264
$ws.Range("J136").Value = 5253.2
$ws.Range("K136").Value = 10678.5792
$ws.Range("L136").Value = 15759.6
$ws.Range("M136").Value = -8128.5792
$ws.Range("N136").Value = -20859.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 406505.12
$ws.Range("I22").Value = 589.6875
$ws.Range("K22").Value = 589.6875
$ws.Range("M22").Value = -416.6875
$ws.Range("H130").Value = 36500
$ws.Range("J130").Value = 36500
$ws.Range("L130").Value = 36500
$ws.Range("N130").Value = -46540
$ws.Range("H134").Value = 2092.4614
$ws.Range("I134").Value = 1731.6
$ws.Range("K134").Value = 5194.799999999999
$ws.Range("M134").Value = -2659.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 57504.844
$ws.Range("I31").Value = 73173.14
$ws.Range("K31").Value = 73173.14
$ws.Range("M31").Value = -72878.14
$ws.Range("H34").Value = 57504.844
$ws.Range("I34").Value = 73173.14
$ws.Range("K34").Value = 73173.14
$ws.Range("M34").Value = -72971.14
$ws.Range("H39").Value = 24399.4
$ws.Range("J39").Value = 55000
$ws.Range("L39").Value = 55000
$ws.Range("N39").Value = -55782
$ws.Range("H49").Value = 24399.4
$ws.Range("J49").Value = 55000
$ws.Range("L49").Value = 55000
$ws.Range("N49").Value = -55364
$ws.Range("H58").Value = 3151.65
$ws.Range("I58").Value = 3348.2354
$ws.Range("K58").Value = 3348.2354
$ws.Range("M58").Value = -3145.2354
$ws.Range("H99").Value = 6983.048
$ws.Range("I99").Value = 5732.1816
$ws.Range("J99").Value = 8359
$ws.Range("K99").Value = 5732.1816
$ws.Range("L99").Value = 8359
$ws.Range("M99").Value = -4234.1816
$ws.Range("N99").Value = -11355
$ws.Range("H122").Value = 1047.7333
$ws.Range("I122").Value = 979.7143
$ws.Range("K122").Value = 2939.1429
$ws.Range("M122").Value = -489.1428999999998
$ws.Range("H126").Value = 6983.048
$ws.Range("I126").Value = 5732.1816
$ws.Range("J126").Value = 8359
$ws.Range("K126").Value = 17196.5448
$ws.Range("L126").Value = 25077
$ws.Range("M126").Value = -14726.5448
$ws.Range("N126").Value = -30017
$ws.Range("H132").Value = 3611.318
$ws.Range("I132").Value = 3444.6843
$ws.Range("K132").Value = 10334.0529
$ws.Range("M132").Value = -7804.052899999999
$ws.Range("H136").Value = 3151.65
$ws.Range("I136").Value = 3348.2354
$ws.Range("K136").Value = 10044.7062
$ws.Range("M136").Value = -7494.706200000001
$ws.Range("H137").Value = 93484.28999999999
$ws.Range("J137").Value = 93484.28999999999
$ws.Range("L137").Value = 93484.28999999999
$ws.Range("N137").Value = -103684.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 19855600
$ws.Range("J9").Value = 69500
$ws.Range("L9").Value = 208500
$ws.Range("N9").Value = -208948
$ws.Range("H121").Value = 30000
$ws.Range("I121").Value = 30000
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 90000
$ws.Range("L121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -88690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3476.2856
$ws.Range("I132").Value = 3128.5
$ws.Range("J132").Value = 4345.75
$ws.Range("K132").Value = 9385.5
$ws.Range("L132").Value = 13037.25
$ws.Range("M132").Value = -6855.5
$ws.Range("N132").Value = -18097.25
$ws.Range("H136").Value = 43749.875
$ws.Range("J136").Value = 43749.875
$ws.Range("L136").Value = 131249.625
$ws.Range("N136").Value = -136349.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1964.625
$ws.Range("I46").Value = 1314.25
$ws.Range("J46").Value = 2615
$ws.Range("K46").Value = 1314.25
$ws.Range("L46").Value = 2615
$ws.Range("M46").Value = -1126.25
$ws.Range("N46").Value = -2991

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 774.94446
$ws.Range("I113").Value = 741.4
$ws.Range("K113").Value = 2224.2
$ws.Range("M113").Value = -54.19999999999982
$ws.Range("H132").Value = 3368.25
$ws.Range("I132").Value = 3185.5557
$ws.Range("K132").Value = 9556.667099999999
$ws.Range("M132").Value = -7026.667099999999
